$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) "OPO-Leader einbauen" paragraph: bold+underline -> strike-through
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "OPO-Leader einbauen*") {
        $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
               '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
               '<pkg:xmlData>' + `
               '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
               '<w:body>' + `
               '<w:p w:rsidR="00FB2C90" w:rsidRPr="00BE1BB6" w:rsidRDefault="00FB2C90" w:rsidP="00032A5C">' + `
               '<w:pPr><w:rPr><w:strike/><w:sz w:val="20"/><w:lang w:val="en-GB"/></w:rPr></w:pPr>' + `
               '<w:r w:rsidRPr="00BE1BB6"><w:rPr><w:strike/><w:sz w:val="20"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">OPO-Leader </w:t></w:r>' + `
               '<w:proofErr w:type="spellStart"/>' + `
               '<w:r w:rsidRPr="00BE1BB6"><w:rPr><w:strike/><w:sz w:val="20"/><w:lang w:val="en-GB"/></w:rPr><w:t>einbauen</w:t></w:r>' + `
               '<w:proofErr w:type="spellEnd"/>' + `
               '</w:p>' + `
               '</w:body></w:document>' + `
               '</pkg:xmlData></pkg:part></pkg:package>'
        $p.Range.InsertXML($xml)
        break
    }
}

# ---------------------------------------------------------------------------
# 2) Remove the trailing " (Discuss)" run after "... hässlich aus!"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute(" (Discuss)", $false, $false, $false, $false, $false, `
                         $true, 1, $false, "", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) "Leader vorerst in fixe, vordefinierte Slots einordnen" paragraph:
#    bold+underline -> strike-through
# ---------------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Leader vorerst in fixe, vordefinierte Slots einordnen*") {
        $xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
               '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
               '<pkg:xmlData>' + `
               '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
               '<w:body>' + `
               '<w:p w:rsidR="005305B4" w:rsidRPr="00CB62D7" w:rsidRDefault="005305B4" w:rsidP="00032A5C">' + `
               '<w:pPr><w:rPr><w:strike/><w:sz w:val="20"/><w:lang w:val="en-GB"/></w:rPr></w:pPr>' + `
               '<w:r w:rsidRPr="00CB62D7"><w:rPr><w:strike/><w:sz w:val="20"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve">Leader </w:t></w:r>' + `
               '<w:proofErr w:type="spellStart"/>' + `
               '<w:r w:rsidRPr="00CB62D7"><w:rPr><w:strike/><w:sz w:val="20"/><w:lang w:val="en-GB"/></w:rPr><w:t>vorerst</w:t></w:r>' + `
               '<w:proofErr w:type="spellEnd"/>' + `
               '<w:r w:rsidRPr="00CB62D7"><w:rPr><w:strike/><w:sz w:val="20"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> in fixe, </w:t></w:r>' + `
               '<w:proofErr w:type="spellStart"/>' + `
               '<w:r w:rsidRPr="00CB62D7"><w:rPr><w:strike/><w:sz w:val="20"/><w:lang w:val="en-GB"/></w:rPr><w:t>vordefinierte</w:t></w:r>' + `
               '<w:proofErr w:type="spellEnd"/>' + `
               '<w:r w:rsidRPr="00CB62D7"><w:rPr><w:strike/><w:sz w:val="20"/><w:lang w:val="en-GB"/></w:rPr><w:t xml:space="preserve"> Slots </w:t></w:r>' + `
               '<w:proofErr w:type="spellStart"/>' + `
               '<w:r w:rsidRPr="00CB62D7"><w:rPr><w:strike/><w:sz w:val="20"/><w:lang w:val="en-GB"/></w:rPr><w:t>einordnen</w:t></w:r>' + `
               '<w:proofErr w:type="spellEnd"/>' + `
               '</w:p>' + `
               '</w:body></w:document>' + `
               '</pkg:xmlData></pkg:part></pkg:package>'
        $p.Range.InsertXML($xml)
        break
    }
}
